# Add season record columns (Wins, Losses, Ties) to the team statistics sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the new
# header cells so they get the same bold / bordered / centered style, then
# set their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every team/player row (2-39) gets the same season record values.
$ws.Range("AD2:AD39").Value = 91
$ws.Range("AE2:AE39").Value = 71
$ws.Range("AF2:AF39").Value = 0
